$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Shape "Text Placeholder 5" (body bullet list)
$bodyShape = $s.Shapes.Item(2)
$bodyRange = $bodyShape.TextFrame.TextRange

# Insert a new bullet paragraph after "Identify where the functions are being
# called" (paragraph 3), leaving the trailing empty paragraph untouched.
$thirdParagraph = $bodyRange.Paragraphs(3)
$newBullet = [char]13 + "Submit your answers to our Attendance assignment for today" + [char]0x2019 + "s class"
[void]$thirdParagraph.InsertAfter($newBullet)

# Grow the placeholder to fit the extra line of text.
$bodyShape.Height = 128.02217

# Shape "TextBox 2" (code sample) moves down to make room for the new bullet.
$codeShape = $s.Shapes.Item(3)
$codeShape.Top = 189.45256
